# KEY ORGANIZERS.xlsx update
#  - Rename "Tech Lead" role -> "Technical Lead" (OPERATIONS TEAM)
#  - Fill in missing Team/Course Background for Show Caller & Tech Officer rows (OPERATIONS TEAM)
#  - Add new team member Darren Bergado (STUDENT SUCCESS TEAM)
#  - Rename "Manager, Marketing Strategy" -> "Marketing Director" (MARKETING & COMMUNICATIONS TEAM)
#  - Remove Kohulan Thevananthan row (ISSUE MANAGEMENT TEAM)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# OPERATIONS TEAM
# ---------------------------------------------------------------
$opsSheet = $wb.Worksheets.Item("OPERATIONS TEAM")

# Haru Ai Okabe, row 12: Role "Tech Lead" -> "Technical Lead"
$opsSheet.Cells.Item(12, 4).Value = "Technical Lead"

# Soumya Saraswati, row 15: fill Team / Course Background
$opsSheet.Cells.Item(15, 3).Value = "Broadcast"
$opsSheet.Cells.Item(15, 5).Value = "Event and Media Production"

# Sanjan, row 16: fill Team / Course Background
$opsSheet.Cells.Item(16, 3).Value = "Broadcast"
$opsSheet.Cells.Item(16, 5).Value = "Event and Media Production"

# ---------------------------------------------------------------
# STUDENT SUCCESS TEAM - add Darren Bergado as a new row
# ---------------------------------------------------------------
$studentSuccessSheet = $wb.Worksheets.Item("STUDENT SUCCESS TEAM")
$newRow = $studentSuccessSheet.UsedRange.Rows.Count + 1

$studentSuccessSheet.Cells.Item($newRow, 1).Value = "Darren Bergado"
$studentSuccessSheet.Cells.Item($newRow, 2).Value = "dbergad1@my.centennialcollege.ca"
$studentSuccessSheet.Cells.Item($newRow, 3).Value = "Student Success"
$studentSuccessSheet.Cells.Item($newRow, 4).Value = "Student Sucess Officer"
$studentSuccessSheet.Cells.Item($newRow, 5).Value = "Electrical Engineering"

# ---------------------------------------------------------------
# MARKETING & COMMUNICATIONS TEAM
# ---------------------------------------------------------------
$marketingSheet = $wb.Worksheets.Item("MARKETING & COMMUNICATIONS TEAM")

# Vidhi Mehta, row 2: Role "Manager, Marketing Strategy" -> "Marketing Director"
$marketingSheet.Cells.Item(2, 4).Value = "Marketing Director"

# ---------------------------------------------------------------
# ISSUE MANAGEMENT TEAM - remove Kohulan Thevananthan row
# ---------------------------------------------------------------
$issueSheet = $wb.Worksheets.Item("ISSUE MANAGEMENT TEAM")
$issueSheet.Rows.Item(2).Delete()
